$wb = $excel.ActiveWorkbook

# New "interest count" (F column) values for rows 2..13 after the shift.
$newF = @{
    2  = 430
    3  = 1450
    4  = 948
    5  = 60
    6  = 2097
    7  = 34
    8  = 1270
    9  = 62
    10 = 113
    11 = 35
    12 = 305
    13 = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2 (the "宜春·融荟城难忘今宵汉文化节" entry) was removed; everything below
    # shifts up one row, so deleting row 2 reproduces the rest of the table.
    $ws.Rows.Item(2).Delete()

    # Column A is a sequential row index (0,1,2,...), independent of which
    # event occupies the row; restore it after the shift.
    for ($row = 2; $row -le 13; $row++) {
        $ws.Cells.Item($row, 1).Value = $row - 1
    }

    # Refresh the "想去人数" (interest count) values that changed between scrapes.
    foreach ($row in $newF.Keys) {
        $ws.Cells.Item($row, 6).Value = $newF[$row]
    }
}
